$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-03 Friday" "2025-10-04 Saturday"

Replace-Text "325×9=" "549×5="
Replace-Text "208×8=" "694×7="
Replace-Text "867×6=" "217×9="
Replace-Text "967×4=" "249×7="
Replace-Text "545×2=" "511×8="

Replace-Text "158×6=" "478×6="
Replace-Text "726×7=" "386×3="
Replace-Text "276×4=" "280×2="
Replace-Text "392×4=" "703×4="
Replace-Text "268×5=" "922×4="

Replace-Text "677×7=" "113×9="
Replace-Text "602×5=" "197×5="
Replace-Text "142×5=" "422×4="
Replace-Text "652×2=" "949×7="
Replace-Text "815×4=" "264×8="

Replace-Text "517×4=" "783×6="
Replace-Text "535×2=" "381×9="
Replace-Text "494×3=" "648×3="
Replace-Text "143×8=" "779×6="
Replace-Text "750×6=" "812×6="

Replace-Text "346×2=" "813×6="
Replace-Text "536×3=" "270×6="
Replace-Text "380×4=" "395×5="
Replace-Text "345×7=" "420×3="
Replace-Text "516×3=" "336×3="
